# Generate Report for Handoff
# Updates the localization-status report with a new handoff run:
#  - new guid-based file id (9c4ab789-94e7-47ca-a82a-1e2b7f44d711) replacing
#    the previous one (b821fe0e-7084-48b4-99ad-0bb2568f4e23)
#  - refreshed handoff timestamps
#  - refreshed xliff content-hash filenames for zh-cn / de-de

$wb = $excel.ActiveWorkbook

$oldId = "b821fe0e-7084-48b4-99ad-0bb2568f4e23"
$newId = "9c4ab789-94e7-47ca-a82a-1e2b7f44d711"

$oldHash = "7c16ef2c4765438be35b99e50685aaa6da712cd4"
$newHash = "52d7f6b687d20326b5247ff129c98eff03c6e25e"

$newFileName = "$newId.md"
$newPathAndName = "e2e\$newId.md"

$newGenerateDate = "2016-08-13 01:11:26"

$newZhHandoffFile = "$newId.$newHash.zh-cn.xlf"
$newZhHandoffDate = "2016-08-13 01:11:18"

$newDeHandoffFile = "$newId.$newHash.de-de.xlf"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathAndName
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newPathAndName
$wsOverview.Range("G2").Value = $newGenerateDate

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newFileName
$wsZh.Hyperlinks.Item(1).TextToDisplay = $newFileName
$wsZh.Range("G2").Value = $newZhHandoffFile
$wsZh.Range("H2").Value = $newZhHandoffDate

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newFileName
$wsDe.Hyperlinks.Item(1).TextToDisplay = $newFileName
$wsDe.Range("G2").Value = $newDeHandoffFile
